# Auto-generated edit script applying numeric updates to Halicarnassus_Profits workbook
# Source: diff of Sheets/Halicarnassus_Profits.xlsx (multi-sheet workbook: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# Updates quest-profit tracking cells (H,I,J,K,L,M,N columns) across all 8 sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4970.5
$ws.Range("I40").Value = 4348.857
$ws.Range("K40").Value = 4348.857
$ws.Range("M40").Value = -4173.857
$ws.Range("H43").Value = 2431.52
$ws.Range("I43").Value = 2419.1765
$ws.Range("J43").Value = 2457.75
$ws.Range("K43").Value = 2419.1765
$ws.Range("L43").Value = 2457.75
$ws.Range("M43").Value = -2350.1765
$ws.Range("N43").Value = -2595.75
$ws.Range("H74").Value = 10997
$ws.Range("I74").Value = 13163
$ws.Range("K74").Value = 13163
$ws.Range("M74").Value = -12227
$ws.Range("H76").Value = 3897.5
$ws.Range("I76").Value = 3894.5
$ws.Range("J76").Value = 3899
$ws.Range("K76").Value = 3894.5
$ws.Range("L76").Value = 3899
$ws.Range("M76").Value = -3579.5
$ws.Range("N76").Value = -4529
$ws.Range("H77").Value = 10997
$ws.Range("I77").Value = 13163
$ws.Range("K77").Value = 65815
$ws.Range("M77").Value = -61135
$ws.Range("H79").Value = 3897.5
$ws.Range("I79").Value = 3894.5
$ws.Range("J79").Value = 3899
$ws.Range("K79").Value = 3894.5
$ws.Range("L79").Value = 3899
$ws.Range("M79").Value = -2802.5
$ws.Range("N79").Value = -6083

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 12020238
$ws.Range("J6").Value = 11667062
$ws.Range("L6").Value = 11667062
$ws.Range("N6").Value = -11667408
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H44").Value = 12775.6875
$ws.Range("J44").Value = 12775.6875
$ws.Range("L44").Value = 12775.6875
$ws.Range("N44").Value = -13751.6875
$ws.Range("H55").Value = 14168
$ws.Range("H97").Value = 1125.4286
$ws.Range("I97").Value = 1146.3334
$ws.Range("K97").Value = 1146.3334
$ws.Range("M97").Value = -650.3334
$ws.Range("H132").Value = 909.4
$ws.Range("I132").Value = 909.4
$ws.Range("K132").Value = 2728.2
$ws.Range("M132").Value = -198.1999999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 2950
$ws.Range("I54").Value = 2950
$ws.Range("K54").Value = 2950
$ws.Range("M54").Value = -2466
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("N94").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 184.85715
$ws.Range("I7").Value = 123.75
$ws.Range("J7").Value = 266.33334
$ws.Range("K7").Value = 123.75
$ws.Range("L7").Value = 266.33334
$ws.Range("M7").Value = -10.75
$ws.Range("N7").Value = -492.33334
$ws.Range("H62").Value = 5502.5
$ws.Range("I62").Value = 3005
$ws.Range("J62").Value = 8000
$ws.Range("K62").Value = 3005
$ws.Range("L62").Value = 8000
$ws.Range("M62").Value = -2381
$ws.Range("N62").Value = -9248
$ws.Range("H65").Value = 5502.5
$ws.Range("I65").Value = 3005
$ws.Range("J65").Value = 8000
$ws.Range("K65").Value = 15025
$ws.Range("L65").Value = 40000
$ws.Range("M65").Value = -11905
$ws.Range("N65").Value = -46240
$ws.Range("H107").Value = 431.6842
$ws.Range("I107").Value = 376.58823
$ws.Range("K107").Value = 376.58823
$ws.Range("M107").Value = 1543.41177

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("H68").Value = 708.6
$ws.Range("I68").Value = 585.75
$ws.Range("J68").Value = 1200
$ws.Range("K68").Value = 1757.25
$ws.Range("L68").Value = 3600
$ws.Range("M68").Value = -946.25
$ws.Range("N68").Value = -5222
$ws.Range("H71").Value = 708.6
$ws.Range("I71").Value = 585.75
$ws.Range("J71").Value = 1200
$ws.Range("K71").Value = 5271.75
$ws.Range("L71").Value = 10800
$ws.Range("M71").Value = -1215.75
$ws.Range("N71").Value = -18912
$ws.Range("H108").Value = 2374.5
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 2374.5
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 7123.5
$ws.Range("M108").ClearContents()
$ws.Range("N108").Value = -12883.5
$ws.Range("H114").Value = 4575.857
$ws.Range("J114").Value = 4505.1665
$ws.Range("L114").Value = 13515.4995
$ws.Range("N114").Value = -20023.4995
$ws.Range("H121").Value = 400
$ws.Range("I121").Value = 400
$ws.Range("K121").Value = 1200
$ws.Range("M121").Value = 110
$ws.Range("H128").Value = 1979899
$ws.Range("I128").Value = 1979899
$ws.Range("K128").Value = 5939697
$ws.Range("M128").Value = -5934717
$ws.Range("H131").Value = 1144.7273
$ws.Range("J131").Value = 1266.5
$ws.Range("L131").Value = 3799.5
$ws.Range("N131").Value = -13879.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1677.8
$ws.Range("I97").Value = 1598.5
$ws.Range("K97").Value = 1598.5
$ws.Range("M97").Value = -1102.5
$ws.Range("H98").Value = 30000
$ws.Range("J98").Value = 30000
$ws.Range("L98").Value = 30000
$ws.Range("N98").Value = -35990
$ws.Range("H99").Value = 9000
$ws.Range("I99").Value = 8000
$ws.Range("K99").Value = 8000
$ws.Range("M99").Value = -5754
$ws.Range("H104").Value = 35000
$ws.Range("I104").Value = 45000
$ws.Range("J104").Value = 25000
$ws.Range("K104").Value = 45000
$ws.Range("L104").Value = 25000
$ws.Range("M104").Value = -41506
$ws.Range("N104").Value = -31988
$ws.Range("H108").Value = 72499.5
$ws.Range("J108").Value = 72499.5
$ws.Range("L108").Value = 72499.5
$ws.Range("N108").Value = -80179.5
$ws.Range("H112").Value = 20000
$ws.Range("I112").Value = 20000
$ws.Range("K112").Value = 20000
$ws.Range("M112").Value = -18892
$ws.Range("H113").Value = 6518.9
$ws.Range("I113").Value = 1297.25
$ws.Range("J113").Value = 10000
$ws.Range("K113").Value = 1297.25
$ws.Range("L113").Value = 10000
$ws.Range("M113").Value = 872.75
$ws.Range("N113").Value = -14340
$ws.Range("H126").Value = 3419.5
$ws.Range("I126").Value = 3419.5
$ws.Range("K126").Value = 10258.5
$ws.Range("M126").Value = -7788.5
$ws.Range("H132").Value = 42215.07
$ws.Range("I132").Value = 51510.434
$ws.Range("K132").Value = 154531.302
$ws.Range("M132").Value = -152001.302

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9560.385
$ws.Range("I40").Value = 8928.700000000001
$ws.Range("K40").Value = 8928.700000000001
$ws.Range("M40").Value = -8792.700000000001
$ws.Range("H46").Value = 5665
$ws.Range("J46").Value = 5665
$ws.Range("L46").Value = 5665
$ws.Range("N46").Value = -6041
$ws.Range("H56").Value = 32499.5
$ws.Range("J56").Value = 29999
$ws.Range("L56").Value = 29999
$ws.Range("N56").Value = -31381
$ws.Range("H122").Value = 4233.875
$ws.Range("I122").Value = 4233.875
$ws.Range("K122").Value = 12701.625
$ws.Range("M122").Value = -10251.625
$ws.Range("H136").Value = 6507
$ws.Range("I136").Value = 6157.7
$ws.Range("K136").Value = 18473.1
$ws.Range("M136").Value = -15923.1

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 24875
$ws.Range("J48").Value = 29833.334
$ws.Range("L48").Value = 29833.334
$ws.Range("N48").Value = -30971.334
$ws.Range("H107").Value = 1373.5
$ws.Range("I107").Value = 1373.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 4120.5
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -2200.5
$ws.Range("N107").ClearContents()
$ws.Range("H122").Value = 4588.1177
$ws.Range("I122").Value = 3998.8
$ws.Range("J122").Value = 4833.6665
$ws.Range("K122").Value = 11996.4
$ws.Range("L122").Value = 14500.9995
$ws.Range("M122").Value = -9546.400000000001
$ws.Range("N122").Value = -19400.9995

